$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-16 Monday" "2024-12-17 Tuesday"

Replace-Text "433÷6=72, 1" "375÷3=125, 0"
Replace-Text "985÷6=164, 1" "428÷2=214, 0"
Replace-Text "301÷9=33, 4" "839÷7=119, 6"
Replace-Text "971÷7=138, 5" "906÷4=226, 2"
Replace-Text "196÷3=65, 1" "245÷6=40, 5"

Replace-Text "551÷7=78, 5" "602÷5=120, 2"
Replace-Text "639÷8=79, 7" "465÷8=58, 1"
Replace-Text "652÷2=326, 0" "558÷5=111, 3"
Replace-Text "936÷8=117, 0" "665÷2=332, 1"
Replace-Text "262÷3=87, 1" "131÷3=43, 2"

Replace-Text "410÷9=45, 5" "353÷8=44, 1"
Replace-Text "484÷2=242, 0" "747÷2=373, 1"
Replace-Text "234÷8=29, 2" "621÷7=88, 5"
Replace-Text "662÷3=220, 2" "894÷6=149, 0"
Replace-Text "889÷2=444, 1" "235÷4=58, 3"

Replace-Text "385÷2=192, 1" "873÷6=145, 3"
Replace-Text "206÷6=34, 2" "850÷4=212, 2"
Replace-Text "796÷3=265, 1" "605÷8=75, 5"
Replace-Text "663÷8=82, 7" "184÷4=46, 0"
Replace-Text "376÷7=53, 5" "661÷3=220, 1"

Replace-Text "344÷4=86, 0" "802÷4=200, 2"
Replace-Text "621÷2=310, 1" "700÷4=175, 0"
Replace-Text "866÷2=433, 0" "458÷6=76, 2"
Replace-Text "895÷3=298, 1" "826÷4=206, 2"
Replace-Text "907÷3=302, 1" "653÷9=72, 5"

Write-Output "done"
